# Update faturamento_diario_lojas with latest daily/total figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("P2").Value = 13557.03
$ws.Range("AG2").Value = 134011.82

# Row 3 - Bibi Cell Vieiralves
$ws.Range("P3").Value = 10000
$ws.Range("AG3").Value = 56967.11

# Row 4 - Bibi Cell Manauara
$ws.Range("O4").Value = 2912.5
$ws.Range("P4").Value = 2022
$ws.Range("AG4").Value = 46206.9

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("P5").Value = 6521.01
$ws.Range("AG5").Value = 41663.87

# Row 6 - total
$ws.Range("O6").Value = 21313.41
$ws.Range("P6").Value = 32100.04
$ws.Range("AG6").Value = 278849.7
